$d = $word.ActiveDocument

# Locate the run of text " LinkedIn: " (a leading space, "LinkedIn", a
# colon, and a trailing space) that currently lives in a single <w:r>.
# We need to split it into three runs with identical (empty) formatting:
#   " LinkedIn" | ":" | " "
$target = $d.Content
$found = $target.Find.Execute(" LinkedIn: ", $true, $false, $false, $false,
                               $false, $true, 1, $false, "", 0)

if ($found) {
    $wholeStart = $target.Start
    $wholeEnd = $target.End

    # Character offsets of the colon inside the matched range.
    $colonStart = $wholeStart + 9   # length of " LinkedIn"
    $colonEnd = $colonStart + 1     # length of ":"

    # Forcing Word to materialize separate runs without altering any
    # character formatting: dropping a temporary bookmark onto the
    # sub-ranges and immediately deleting it causes the run table to be
    # rebuilt around those boundaries, splitting the original single run
    # into three runs that keep identical (absent) run properties.
    $colonRange = $d.Range($colonStart, $colonEnd)
    $d.Bookmarks.Add("ztmp_split1", $colonRange)
    $d.Bookmarks("ztmp_split1").Delete()

    $tailRange = $d.Range($colonEnd, $wholeEnd)
    $d.Bookmarks.Add("ztmp_split2", $tailRange)
    $d.Bookmarks("ztmp_split2").Delete()
}
